$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-CellText 'D2' '63.763.76'
Set-CellText 'E2' '  -0.06%  '
Set-CellText 'D3' '2.621.04'
Set-CellText 'E3' '  -0.10%  '
Set-CellText 'E4' '  +0.03%  '
Set-CellText 'D5' '595.09'
Set-CellText 'E5' '  -0.41%  '
Set-CellText 'D6' '151.01'
Set-CellText 'E6' '  +0.48%  '
Set-CellText 'E7' '  +0.02%  '
Set-CellText 'E8' '  -0.28%  '
Set-CellText 'E9' '  +4.29%  '
Set-CellText 'E10' '  +3.74%  '
Set-CellText 'E11' '  +2.46%  '
Set-CellText 'E12' '  +0.95%  '
Set-CellText 'D13' '27.82'
Set-CellText 'E13' '  +0.45%  '
Set-CellText 'D14' '3.093.85'
Set-CellText 'E14' '  -0.04%  '
Set-CellText 'D15' '63.682.52'
Set-CellText 'E15' '  +0.06%  '
Set-CellText 'D16' '0.0000169'
Set-CellText 'E16' '  +13.62%  '
Set-CellText 'D17' '2.650.49'
Set-CellText 'E17' '  +1.06%  '
Set-CellText 'D18' '12.16'
Set-CellText 'E18' '  -1.05%  '
Set-CellText 'E19' '  +3.32%  '
Set-CellText 'D20' '346.76'
Set-CellText 'E20' '  -0.92%  '
Set-CellText 'D21' '6.99'
Set-CellText 'E21' '  +1.51%  '
Set-CellText 'D22' '0.999'
Set-CellText 'E22' '  +0.16%  '
Set-CellText 'D23' '67.27'
Set-CellText 'E23' '  +1.48%  '
Set-CellText 'E24' '  -3.25%  '
Set-CellText 'E25' '  +0.22%  '
Set-CellText 'D26' '9.14'
Set-CellText 'E26' '  -0.46%  '
Set-CellText 'D27' '8.27'
Set-CellText 'E27' '  +0.90%  '
Set-CellText 'D28' '549.60'
Set-CellText 'E28' '  -2.47%  '
Set-CellText 'D29' '0.162'
Set-CellText 'E29' '  -1.55%  '
Set-CellText 'D30' '0.999'
Set-CellText 'E30' '  -0.08%  '
Set-CellText 'D31' '0.0₃0899'
Set-CellText 'E31' '  +6.85%  '
Set-CellText 'D33' '1.81'
Set-CellText 'E33' '  +4.65%  '
Set-CellText 'D34' '5.35'
Set-CellText 'E34' '  +2.66%  '
Set-CellText 'D35' '6.08'
Set-CellText 'E35' '  +0.27%  '
Set-CellText 'E36' '  +2.02%  '
Set-CellText 'D37' '164.07'
Set-CellText 'E37' '  -3.04%  '
Set-CellText 'D38' '19.95'
Set-CellText 'E38' '  +3.14%  '
Set-CellText 'B39' 'FirstDigitalUSD'
Set-CellText 'C39' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText 'D39' '1.00'
Set-CellText 'E39' '  +0.09%  '
Set-CellText 'B40' 'Stacks'
Set-CellText 'C40' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText 'D40' '1.97'
Set-CellText 'E40' '  +1.51%  '
Set-CellText 'E41' '  -0.03%  '
Set-CellText 'D42' '167.48'
Set-CellText 'E42' '  -1.66%  '
Set-CellText 'D43' '4.09'
Set-CellText 'E43' '  +4.66%  '
Set-CellText 'D44' '23.27'
Set-CellText 'E44' '  +8.95%  '
Set-CellText 'E45' '  -2.52%  '
Set-CellText 'D46' '2.18'
Set-CellText 'E46' '  +10.37%  '
Set-CellText 'D47' '0.637'
Set-CellText 'E47' '  +1.09%  '
Set-CellText 'D48' '0.0251'
Set-CellText 'E48' '  +1.61%  '
Set-CellText 'E49' '  +0.08%  '
Set-CellText 'D50' '19.24'
Set-CellText 'E50' '  +0.12%  '
Set-CellText 'D51' '0.0₆0233'
Set-CellText 'E51' '  +19.30%  '
